$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold, bordered, centered) onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header row values
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 8
